# Unit 49 vocabulary additions (30 new rows: 1442-1471).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Phase 1: column A (word IDs) for all 30 new rows first, matching the
# original authoring order reflected in the shared-string table ordering.
$ws.Range("A1442").Value = 'U49_01'
$ws.Range("A1443").Value = 'U49_02'
$ws.Range("A1444").Value = 'U49_03'
$ws.Range("A1445").Value = 'U49_04'
$ws.Range("A1446").Value = 'U49_05'
$ws.Range("A1447").Value = 'U49_06'
$ws.Range("A1448").Value = 'U49_07'
$ws.Range("A1449").Value = 'U49_08'
$ws.Range("A1450").Value = 'U49_09'
$ws.Range("A1451").Value = 'U49_10'
$ws.Range("A1452").Value = 'U49_11'
$ws.Range("A1453").Value = 'U49_12'
$ws.Range("A1454").Value = 'U49_13'
$ws.Range("A1455").Value = 'U49_14'
$ws.Range("A1456").Value = 'U49_15'
$ws.Range("A1457").Value = 'U49_16'
$ws.Range("A1458").Value = 'U49_17'
$ws.Range("A1459").Value = 'U49_18'
$ws.Range("A1460").Value = 'U49_19'
$ws.Range("A1461").Value = 'U49_20'
$ws.Range("A1462").Value = 'U49_21'
$ws.Range("A1463").Value = 'U49_22'
$ws.Range("A1464").Value = 'U49_23'
$ws.Range("A1465").Value = 'U49_24'
$ws.Range("A1466").Value = 'U49_25'
$ws.Range("A1467").Value = 'U49_26'
$ws.Range("A1468").Value = 'U49_27'
$ws.Range("A1469").Value = 'U49_28'
$ws.Range("A1470").Value = 'U49_29'
$ws.Range("A1471").Value = 'U49_30'

# Phase 2: remaining columns per row (B numeric; C/D/E/F text in original
# write order; G reuses an existing part-of-speech string).

$ws.Range("B1442").Value = 49
$ws.Range("C1442").Value = 'Nhà tù'
$ws.Range("D1442").Value = 'Prison'
$ws.Range("E1442").Value = 'He will be in prison for murder'
$ws.Range("F1442").Value = 'in prison for doing something / phải vào tù vì tội gì đó'
$ws.Range("G1442").Value = 'N'

$ws.Range("B1443").Value = 49
$ws.Range("C1443").Value = 'Trốn thoát, thoát khỏi'
$ws.Range("D1443").Value = 'Escape'
$ws.Range("E1443").Value = 'She must escape from him'
$ws.Range("F1443").Value = 'escape from somebody or something'
$ws.Range("G1443").Value = 'V'

$ws.Range("B1444").Value = 49
$ws.Range("C1444").Value = 'Cướp'
$ws.Range("D1444").Value = 'Rob'
$ws.Range("F1444").Value = 'rob somebody or something'
$ws.Range("E1444").Value = 'It is wrong to rob people / trộm cướp tài sản là điều sai trái'
$ws.Range("G1444").Value = 'V'

$ws.Range("B1445").Value = 49
$ws.Range("C1445").Value = 'Đột nhập vào '
$ws.Range("D1445").Value = 'Break into'
$ws.Range("E1445").Value = 'Someone broke into my house'
$ws.Range("F1445").Value = 'break into something / đột nhập vào một cái gì đó'
$ws.Range("G1445").Value = 'V'

$ws.Range("B1446").Value = 49
$ws.Range("C1446").Value = 'Đuổi theo'
$ws.Range("D1446").Value = 'Chase'
$ws.Range("E1446").Value = 'Cats chase mice (chuột)'
$ws.Range("F1446").Value = 'chase somebody or something / đuổi theo ai đó hoặc cái gì đó'
$ws.Range("G1446").Value = 'V'

$ws.Range("B1447").Value = 49
$ws.Range("C1447").Value = 'Che giấu'
$ws.Range("D1447").Value = 'Hide'
$ws.Range("F1447").Value = 'hide somebody or something'
$ws.Range("E1447").Value = 'They hide their money throughout the house (khắp nhà)'
$ws.Range("G1447").Value = 'V'

$ws.Range("B1448").Value = 49
$ws.Range("C1448").Value = 'Tội phạm'
$ws.Range("D1448").Value = 'Criminal'
$ws.Range("E1448").Value = 'The convicted criminal pleaded guilty (nhận tội)'
$ws.Range("F1448").Value = 'Convicted criminal / tội phạm bị kết án'
$ws.Range("G1448").Value = 'N'

$ws.Range("B1449").Value = 49
$ws.Range("C1449").Value = 'Nạn nhân'
$ws.Range("D1449").Value = 'Victim'
$ws.Range("E1449").Value = 'He is a victim of his own greed (Anh ta là nạn nhân của lòng tham của chính mình)'
$ws.Range("F1449").Value = 'a victim of something / nạn nhân của một cái gì đó'
$ws.Range("G1449").Value = 'N'

$ws.Range("B1450").Value = 49
$ws.Range("C1450").Value = 'Luật sư'
$ws.Range("D1450").Value = 'Lawyer'
$ws.Range("E1450").Value = 'He needs to hire a criminal lawyer'
$ws.Range("F1450").Value = 'a criminal lawyer / một luật sư hình sự'
$ws.Range("G1450").Value = 'N'

$ws.Range("B1451").Value = 49
$ws.Range("C1451").Value = 'Có tội, có lỗi'
$ws.Range("D1451").Value = 'Guilty'
$ws.Range("E1451").Value = 'I feel guilty about lying to you'
$ws.Range("F1451").Value = 'guilty about doing something / có lỗi vì làm điều gì đó'
$ws.Range("G1451").Value = 'Adj'

$ws.Range("B1452").Value = 49
$ws.Range("C1452").Value = 'Tòa án'
$ws.Range("D1452").Value = 'Court'
$ws.Range("E1452").Value = 'He can be a witness in a civil court case'
$ws.Range("F1452").Value = 'a civil court case / tòa án dân sự'
$ws.Range("G1452").Value = 'N'

$ws.Range("B1453").Value = 49
$ws.Range("C1453").Value = 'Trường hợp'
$ws.Range("D1453").Value = 'Case'
$ws.Range("E1453").Value = 'In case of emergency, here is my number'
$ws.Range("F1453").Value = 'in case of something / trong trường hợp nào đó'
$ws.Range("G1453").Value = 'N'

$ws.Range("B1454").Value = 49
$ws.Range("C1454").Value = 'Bằng chứng'
$ws.Range("D1454").Value = 'Evidence'
$ws.Range("E1454").Value = 'Scientists are looking for evidence of life on other planets (hành tinh)'
$ws.Range("F1454").Value = 'evidence of something / bằng chứng của một cái gì đó'
$ws.Range("G1454").Value = 'N'

$ws.Range("B1455").Value = 49
$ws.Range("C1455").Value = 'Sự có mặt, sự hiện diện'
$ws.Range("D1455").Value = 'Presence'
$ws.Range("E1455").Value = 'In the presence of her, I feel safe'
$ws.Range("F1455").Value = 'in the presence of something or somebody'
$ws.Range("G1455").Value = 'N'

$ws.Range("B1456").Value = 49
$ws.Range("C1456").Value = 'Tìm kiếm'
$ws.Range("D1456").Value = 'Seek'
$ws.Range("E1456").Value = 'Make sure you seek help if you feel overwhelmed (quá tải)'
$ws.Range("F1456").Value = 'Seek something or somebody / tìm kiếm một cái gì đó'
$ws.Range("G1456").Value = 'V'

$ws.Range("B1457").Value = 49
$ws.Range("C1457").Value = 'An ninh'
$ws.Range("D1457").Value = 'Security'
$ws.Range("E1457").Value = 'National security is a duty of government'
$ws.Range("F1457").Value = 'Nation security / an ninh quốc gia'
$ws.Range("G1457").Value = 'N'

$ws.Range("B1458").Value = 49
$ws.Range("C1458").Value = 'Độc ác, tàn nhẫn'
$ws.Range("D1458").Value = 'Cruel'
$ws.Range("E1458").Value = 'I am sorry that I was cruel to you'
$ws.Range("F1458").Value = 'to be cruel to somebody or something'
$ws.Range("G1458").Value = 'Adj'

$ws.Range("B1459").Value = 49
$ws.Range("C1459").Value = 'Sự tự do'
$ws.Range("D1459").Value = 'Liberty'
$ws.Range("E1459").Value = 'Citizens have the liberty to voice concerns'
$ws.Range("F1459").Value = 'Liberty to do something / tự do làm cái gì đó'
$ws.Range("G1459").Value = 'N'

$ws.Range("B1460").Value = 49
$ws.Range("C1460").Value = 'Tấn công'
$ws.Range("D1460").Value = 'Attack'
$ws.Range("F1460").Value = 'attack somebody / tấn công ai đó'
$ws.Range("E1460").Value = 'The burglar (tên trộm) attacked the civilian (người dân)'
$ws.Range("G1460").Value = 'V'

$ws.Range("B1461").Value = 49
$ws.Range("C1461").Value = 'bạo lực'
$ws.Range("D1461").Value = 'Violent'
$ws.Range("E1461").Value = 'We don''t have to be violent towards one another'
$ws.Range("F1461").Value = 'violent towards or to somebody'
$ws.Range("G1461").Value = 'Adj'

$ws.Range("B1462").Value = 49
$ws.Range("C1462").Value = 'Lan truyền'
$ws.Range("D1462").Value = 'Spread'
$ws.Range("E1462").Value = 'He spreads the word about her promotion'
$ws.Range("F1462").Value = 'spread the word / loan tin'
$ws.Range("G1462").Value = 'V'

$ws.Range("B1463").Value = 49
$ws.Range("C1463").Value = 'Lạm dụng, lợi dụng'
$ws.Range("D1463").Value = 'Abuse'
$ws.Range("E1463").Value = 'Be careful not to abuse alcohol'
$ws.Range("F1463").Value = 'abuse alcohol / lạm dụng rượu'
$ws.Range("G1463").Value = 'V'

$ws.Range("B1464").Value = 49
$ws.Range("C1464").Value = 'Chính phủ'
$ws.Range("D1464").Value = 'Government'
$ws.Range("E1464").Value = 'There will be a high-level meeting among government officials'
$ws.Range("F1464").Value = 'a government official / quan chức chính phủ'
$ws.Range("G1464").Value = 'N'

$ws.Range("B1465").Value = 49
$ws.Range("C1465").Value = 'Luật, quy định'
$ws.Range("E1465").Value = 'What is the law on drinking age?'
$ws.Range("F1465").Value = 'law on something'
$ws.Range("D1465").Value = 'Law'
$ws.Range("G1465").Value = 'N'

$ws.Range("B1466").Value = 49
$ws.Range("C1466").Value = 'Nhanh chóng'
$ws.Range("D1466").Value = 'Rapid'
$ws.Range("E1466").Value = 'The company is experiencing a rapid growth'
$ws.Range("F1466").Value = 'rapid growth / sự tăng trưởng nhanh chóng'
$ws.Range("G1466").Value = 'Adj'

$ws.Range("B1467").Value = 49
$ws.Range("C1467").Value = 'Theo pháp luật'
$ws.Range("D1467").Value = 'Legal'
$ws.Range("E1467").Value = 'We rely on (dựa vào) the legal system for justice'
$ws.Range("F1467").Value = 'the legal system / hệ thống luật pháp'
$ws.Range("G1467").Value = 'Adj'

$ws.Range("B1468").Value = 49
$ws.Range("C1468").Value = 'Sau cùng'
$ws.Range("D1468").Value = 'Ultimate'
$ws.Range("E1468").Value = 'Having a meaningful life is the ultimate target'
$ws.Range("F1468").Value = 'ultimate target / mục tiêu cuối cùng'
$ws.Range("G1468").Value = 'Adj'

$ws.Range("B1469").Value = 49
$ws.Range("C1469").Value = 'Cảnh nghèo nàn'
$ws.Range("D1469").Value = 'Poverty'
$ws.Range("E1469").Value = 'Some countries have extreme poverty'
$ws.Range("F1469").Value = 'extreme poverty / tình trạng nghèo đói cùng cực'
$ws.Range("G1469").Value = 'N'

$ws.Range("B1470").Value = 49
$ws.Range("C1470").Value = 'Ly dị'
$ws.Range("D1470").Value = 'Divorce'
$ws.Range("E1470").Value = 'In the past, wives of kings could never divorce their husbands'
$ws.Range("F1470").Value = 'divorce someone'
$ws.Range("G1470").Value = 'V'

$ws.Range("B1471").Value = 49
$ws.Range("C1471").Value = 'Tội giết người'
$ws.Range("D1471").Value = 'Murder'
$ws.Range("E1471").Value = 'The detective found out who murdered the victim'
$ws.Range("F1471").Value = 'murder somebody / giết ai đó'
$ws.Range("G1471").Value = 'V'

$ws.Range("C1472").Select()
